$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1733333333333333
$ws.Range("C2").Value = 0.6044444444444445
$ws.Range("J2").Value = 0.004444444444444444
$ws.Range("P2").Value = 0.1377777777777778
$ws.Range("S2").Value = 0.08
$ws.Range("B3").Value = 0.02083333333333333
$ws.Range("C3").Value = 0.04861111111111111
$ws.Range("J3").Value = 0.006944444444444444
$ws.Range("P3").Value = 0.7222222222222222
$ws.Range("S3").Value = 0.2013888888888889
$ws.Range("B6").Value = 0.03286384976525822
$ws.Range("D6").Value = 0.0187793427230047
$ws.Range("E6").Value = 0.009389671361502348
$ws.Range("F6").Value = 0.09389671361502347
$ws.Range("J6").Value = 0.215962441314554
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.1690140845070423
$ws.Range("R6").Value = 0.0892018779342723
$ws.Range("S6").Value = 0.3568075117370892
$ws.Range("B7").Value = 0.07929515418502203
$ws.Range("D7").Value = 0.004405286343612335
$ws.Range("F7").Value = 0.03964757709251102
$ws.Range("J7").Value = 0.13215859030837
$ws.Range("O7").Value = 0.00881057268722467
$ws.Range("Q7").Value = 0.1938325991189427
$ws.Range("R7").Value = 0.07488986784140969
$ws.Range("S7").Value = 0.4669603524229075
$ws.Range("B8").Value = 0.07441860465116279
$ws.Range("D8").Value = 0.01627906976744186
$ws.Range("F8").Value = 0.04418604651162791
$ws.Range("J8").Value = 0.1209302325581395
$ws.Range("O8").Value = 0.02093023255813953
$ws.Range("Q8").Value = 0.1604651162790698
$ws.Range("R8").Value = 0.113953488372093
$ws.Range("S8").Value = 0.4488372093023256
$ws.Range("B9").Value = 0.068
$ws.Range("D9").Value = 0.012
$ws.Range("F9").Value = 0.056
$ws.Range("J9").Value = 0.136
$ws.Range("O9").Value = 0.024
$ws.Range("Q9").Value = 0.212
$ws.Range("R9").Value = 0.096
$ws.Range("S9").Value = 0.396
$ws.Range("B10").Value = 0.08061002178649238
$ws.Range("D10").Value = 0.02687000726216413
$ws.Range("F10").Value = 0.06971677559912855
$ws.Range("J10").Value = 0.1350762527233116
$ws.Range("O10").Value = 0.01016702977487291
$ws.Range("Q10").Value = 0.2142338416848221
$ws.Range("R10").Value = 0.1002178649237473
$ws.Range("S10").Value = 0.3631082062454611
$ws.Range("G11").Value = 0.1661538461538462
$ws.Range("J11").Value = 0.04307692307692308
$ws.Range("K11").Value = 0.2
$ws.Range("L11").Value = 0.5723076923076923
$ws.Range("S11").Value = 0.01846153846153846
$ws.Range("G12").Value = 0.7384615384615385
$ws.Range("J12").Value = 0.1897435897435897
$ws.Range("K12").Value = 0.005128205128205128
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.04102564102564103
$ws.Range("G13").Value = 0.6595744680851063
$ws.Range("J13").Value = 0.2978723404255319
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.01630434782608696
$ws.Range("H15").Value = 0.1521739130434783
$ws.Range("I15").Value = 0.09239130434782608
$ws.Range("J15").Value = 0.3641304347826087
$ws.Range("K15").Value = 0.04347826086956522
$ws.Range("M15").Value = 0.0108695652173913
$ws.Range("O15").Value = 0.04347826086956522
$ws.Range("S15").Value = 0.2771739130434783
$ws.Range("F16").Value = 0.01796407185628742
$ws.Range("H16").Value = 0.1856287425149701
$ws.Range("I16").Value = 0.1197604790419162
$ws.Range("J16").Value = 0.4131736526946108
$ws.Range("K16").Value = 0.125748502994012
$ws.Range("M16").Value = 0.02994011976047904
$ws.Range("O16").Value = 0.02395209580838323
$ws.Range("S16").Value = 0.08383233532934131
$ws.Range("F17").Value = 0.006109979633401222
$ws.Range("H17").Value = 0.1588594704684318
$ws.Range("I17").Value = 0.120162932790224
$ws.Range("J17").Value = 0.4276985743380856
$ws.Range("K17").Value = 0.1120162932790224
$ws.Range("M17").Value = 0.01221995926680244
$ws.Range("O17").Value = 0.05091649694501019
$ws.Range("S17").Value = 0.1120162932790224
$ws.Range("F18").Value = 0.01209677419354839
$ws.Range("H18").Value = 0.1975806451612903
$ws.Range("I18").Value = 0.1088709677419355
$ws.Range("J18").Value = 0.4153225806451613
$ws.Range("K18").Value = 0.1129032258064516
$ws.Range("M18").Value = 0.01612903225806452
$ws.Range("N18").Value = 0.004032258064516129
$ws.Range("O18").Value = 0.04032258064516129
$ws.Range("S18").Value = 0.09274193548387097
$ws.Range("F19").Value = 0.01630434782608696
$ws.Range("H19").Value = 0.1894409937888199
$ws.Range("I19").Value = 0.09782608695652174
$ws.Range("J19").Value = 0.4083850931677019
$ws.Range("K19").Value = 0.1125776397515528
$ws.Range("M19").Value = 0.02329192546583851
$ws.Range("N19").Value = 0.002329192546583851
$ws.Range("O19").Value = 0.05667701863354038
$ws.Range("S19").Value = 0.09316770186335403
